{"js": "// The document contains a single 20-row x 5-column table of simple\n// arithmetic equations (e.g. \"43-31=12\"). The edit replaces the text of\n// 100 specific cells with new equations, in row-major (document) order,\n// leaving every other part of the document (title paragraph, table\n// structure, run/paragraph formatting) untouched.\n//\n// newValues[r][c] is the replacement text for the cell at row r, col c\n// (0-based), matching the order the equations appear in the document.\nconst newValues = [\n  [\"80-27=53\", \"78-0=78\", \"95-40=55\", \"34-24=10\", \"14+72=86\"],\n  [\"96-31=65\", \"42-31=11\", \"21+5=26\", \"35+44=79\", \"23-15=8\"],\n  [\"88-9=79\", \"88-14=74\", \"11+75=86\", \"57+16=73\", \"54-43=11\"],\n  [\"26+31=57\", \"85-7=78\", \"7+9=16\", \"57-0=57\", \"13+14=27\"],\n  [\"33+26=59\", \"76-61=15\", \"44+16=60\", \"46+4=50\", \"25+10=35\"],\n  [\"55+0=55\", \"60-31=29\", \"35-29=6\", \"65-36=29\", \"5+16=21\"],\n  [\"10+53=63\", \"86-64=22\", \"54-24=30\", \"77-21=56\", \"28-6=22\"],\n  [\"51-37=14\", \"98-19=79\", \"64+1=65\", \"11+51=62\", \"25-23=2\"],\n  [\"63+30=93\", \"27-24=3\", \"58-15=43\", \"25+5=30\", \"26+69=95\"],\n  [\"43+30=73\", \"0+30=30\", \"94-1=93\", \"17+50=67\", \"54-6=48\"],\n  [\"95-4=91\", \"45+36=81\", \"70-51=19\", \"40-1=39\", \"15+31=46\"],\n  [\"21+73=94\", \"57-23=34\", \"12+34=46\", \"30+46=76\", \"71-28=43\"],\n  [\"25-14=11\", \"29+22=51\", \"33+50=83\", \"78-22=56\", \"71+26=97\"],\n  [\"53+14=67\", \"25-9=16\", \"29+5=34\", \"98-61=37\", \"6+11=17\"],\n  [\"0+48=48\", \"46-15=31\", \"85-47=38\", \"13+25=38\", \"35+22=57\"],\n  [\"26+57=83\", \"78-1=77\", \"12+53=65\", \"86-28=58\", \"25+2=27\"],\n  [\"21+57=78\", \"7+44=51\", \"84-25=59\", \"66-11=55\", \"16+77=93\"],\n  [\"96-46=50\", \"37-23=14\", \"7+67=74\", \"77-77=0\", \"25+2=27\"],\n  [\"33+36=69\", \"96-17=79\", \"16+13=29\", \"34+19=53\", \"91-51=40\"],\n  [\"57+27=84\", \"55-49=6\", \"28+29=57\", \"4+22=26\", \"79-34=45\"],\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nfor (let r = 0; r < newValues.length; r++) {\n  for (let c = 0; c < newValues[r].length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[r][c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document contains a single 20-row x 5-column table of simple\n# arithmetic equations (e.g. \"43-31=12\"). The edit replaces the text of\n# 100 specific cells with new equations, in row-major (document) order,\n# leaving every other part of the document (title paragraph, table\n# structure, run/paragraph formatting) untouched.\n#\n# $newValues[r][c] is the replacement text for the cell at row r, col c\n# (0-based in the array; COM Cell() is 1-based), matching the order the\n# equations appear in the document.\n$newValues = @(\n  @(\"80-27=53\", \"78-0=78\", \"95-40=55\", \"34-24=10\", \"14+72=86\"),\n  @(\"96-31=65\", \"42-31=11\", \"21+5=26\", \"35+44=79\", \"23-15=8\"),\n  @(\"88-9=79\", \"88-14=74\", \"11+75=86\", \"57+16=73\", \"54-43=11\"),\n  @(\"26+31=57\", \"85-7=78\", \"7+9=16\", \"57-0=57\", \"13+14=27\"),\n  @(\"33+26=59\", \"76-61=15\", \"44+16=60\", \"46+4=50\", \"25+10=35\"),\n  @(\"55+0=55\", \"60-31=29\", \"35-29=6\", \"65-36=29\", \"5+16=21\"),\n  @(\"10+53=63\", \"86-64=22\", \"54-24=30\", \"77-21=56\", \"28-6=22\"),\n  @(\"51-37=14\", \"98-19=79\", \"64+1=65\", \"11+51=62\", \"25-23=2\"),\n  @(\"63+30=93\", \"27-24=3\", \"58-15=43\", \"25+5=30\", \"26+69=95\"),\n  @(\"43+30=73\", \"0+30=30\", \"94-1=93\", \"17+50=67\", \"54-6=48\"),\n  @(\"95-4=91\", \"45+36=81\", \"70-51=19\", \"40-1=39\", \"15+31=46\"),\n  @(\"21+73=94\", \"57-23=34\", \"12+34=46\", \"30+46=76\", \"71-28=43\"),\n  @(\"25-14=11\", \"29+22=51\", \"33+50=83\", \"78-22=56\", \"71+26=97\"),\n  @(\"53+14=67\", \"25-9=16\", \"29+5=34\", \"98-61=37\", \"6+11=17\"),\n  @(\"0+48=48\", \"46-15=31\", \"85-47=38\", \"13+25=38\", \"35+22=57\"),\n  @(\"26+57=83\", \"78-1=77\", \"12+53=65\", \"86-28=58\", \"25+2=27\"),\n  @(\"21+57=78\", \"7+44=51\", \"84-25=59\", \"66-11=55\", \"16+77=93\"),\n  @(\"96-46=50\", \"37-23=14\", \"7+67=74\", \"77-77=0\", \"25+2=27\"),\n  @(\"33+36=69\", \"96-17=79\", \"16+13=29\", \"34+19=53\", \"91-51=40\"),\n  @(\"57+27=84\", \"55-49=6\", \"28+29=57\", \"4+22=26\", \"79-34=45\")\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\nfor ($r = 0; $r -lt $newValues.Count; $r++) {\n  $row = $newValues[$r]\n  for ($c = 0; $c -lt $row.Count; $c++) {\n    $cell = $t.Cell($r + 1, $c + 1)\n    $cell.Range.Text = $row[$c]\n  }\n}\n"}
